$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "latest" row (47) reverts to the standard date-time format
$ws.Range("A47").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 48
$ws.Range("A48").NumberFormat = "YYYY-MM-DD"
$ws.Range("A48").Value = 45633
$ws.Range("B48").Value = 122
$ws.Range("C48").Value = 108
$ws.Range("D48").Value = 116
